# Update the "Förändrad" (last-changed) date in column C for every data row,
# and append the display-text argument to the HYPERLINK() formulas in the
# link columns (S,T,U,V,W,X,Y) for the rows that actually have such links.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow  = 224
$newChangedDate = 45186

$linkCols = @("S", "T", "U", "V", "W", "X", "Y")

for ($row = $firstDataRow; $row -le $lastDataRow; $row++) {

    # Column C: bump the "changed" serial date on every row.
    $ws.Range("C" + $row).Value = $newChangedDate

    # Designation code in column A is used as the HYPERLINK friendly name.
    $aCell = $ws.Range("A" + $row)
    $aVal = $aCell.Value()

    foreach ($col in $linkCols) {
        $cell = $ws.Range($col + $row)
        if ($cell.HasFormula) {
            $f = $cell.Formula
            # Only touch plain single-argument HYPERLINK(...) formulas that
            # don't already carry a friendly-name second argument.
            if ($f.StartsWith("=HYPERLINK(") -and -not $f.Contains(",")) {
                $newFormula = $f.Substring(0, $f.Length - 1) + ', "' + $aVal + '")'
                $cell.Formula = $newFormula
            }
        }
    }
}
